$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the test status for the Step2Fields row (row 4): it finished testing
# and moved from "Testing" to "Automated".
$ws.Range("B4").Value = 8
$ws.Range("D4").Value = "Automated"

# Move the active selection to D5, as recorded after finishing the edit.
$ws.Range("D5").Select()
